$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.850.52"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.376.85"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.71"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.05"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.377.37"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.945.89"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.95"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.368.94"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.932.19"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.65"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.22"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.74"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.508.00"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.546"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  +10.17%  "
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.24"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.53"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.75"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.70"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0757"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.771"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.70"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.96"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.31"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.532.03"
$ws.Range("E48").Value = "  +7.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.38"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  +3.00%  "
